$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.272.28"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "1.662.34"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  +0.76%  "
$ws.Range("D5").Value = "218.29"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").Value = "0.5325"
$ws.Range("E6").Value = "  +0.89%  "
$ws.Range("D8").Value = "0.2634"
$ws.Range("E8").Value = "  +1.02%  "
$ws.Range("D9").Value = "0.06355"
$ws.Range("D10").Value = "20.50"
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("D11").Value = "0.07824"
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("D12").Value = "4.563"
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("D13").Value = "1.698.13"
$ws.Range("E13").Value = "  +3.67%  "
$ws.Range("D14").Value = "1.889.94"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("D15").Value = "0.5532"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").Value = "0.0₅8177"
$ws.Range("E16").Value = "  -0.76%  "
$ws.Range("D17").Value = "65.64"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("D19").Value = "4.675"
$ws.Range("E19").Value = "  +2.02%  "
$ws.Range("D20").Value = "193.00"
$ws.Range("E20").Value = "  +0.40%  "
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("D22").Value = "6.019"
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("E23").Value = "  +0.71%  "
$ws.Range("D24").Value = "145.80"
$ws.Range("E24").Value = "  +2.77%  "
$ws.Range("D25").Value = "0.1223"
$ws.Range("E25").Value = "  -2.29%  "
$ws.Range("D26").Value = "7.181"
$ws.Range("E26").Value = "  -1.39%  "
$ws.Range("E27").Value = "  -0.84%  "
$ws.Range("D28").Value = "1.483"
$ws.Range("E28").Value = "  +3.01%  "
$ws.Range("D29").Value = "0.05886"
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("D31").Value = "3.585"
$ws.Range("E31").Value = "  +1.72%  "
$ws.Range("D32").Value = "3.270"
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("E33").Value = "  +1.35%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "0.9600"
$ws.Range("E34").Value = "  +0.56%  "
$ws.Range("B35").Value = "MXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D35").Value = "2.822"
$ws.Range("E35").Value = "  +1.07%  "
$ws.Range("E36").Value = "  +0.48%  "
$ws.Range("D37").Value = "0.5786"
$ws.Range("E37").Value = "  +1.22%  "
$ws.Range("E38").Value = "  -1.37%  "
$ws.Range("D39").Value = "0.8619"
$ws.Range("E39").Value = "  +1.46%  "
$ws.Range("D40").Value = "5.831"
$ws.Range("E40").Value = "  +0.82%  "
$ws.Range("E41").Value = "  +0.70%  "
$ws.Range("D42").Value = "1.044.92"
$ws.Range("E42").Value = "  +1.62%  "
$ws.Range("D43").Value = "104.12"
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("D44").Value = "1.801.57"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "57.44"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").Value = "1.010"
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₈105"
$ws.Range("E47").Value = "  -5.71%  "
$ws.Range("E48").Value = "  +1.87%  "
$ws.Range("D49").Value = "8.011"
$ws.Range("E49").Value = "  +1.87%  "
$ws.Range("D50").Value = "0.05162"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").Value = "1.429"
$ws.Range("E51").Value = "  -3.69%  "
